$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Title "Test Plan" paragraph: remove centered alignment
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.ParagraphFormat.Alignment = 0

# ---------------------------------------------------------------------------
# 2. RF14 -> RF8 text change + move the _GoBack bookmark
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "RF14: Visualizzazione del curriculum di un candidato",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "RF8: Visualizzazione delle candidature effettuate", 2) | Out-Null

# Remove the old _GoBack bookmark (currently on the RF16 paragraph)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Re-create _GoBack at the end of the RF8 paragraph (empty range == collapsed bookmark)
$found = $d.Content
$found.Find.Execute("RF8: Visualizzazione delle candidature effettuate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRange = $d.Range($found.End, $found.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

Write-Host "stage1 done"
